# consulta_cnpj_esalva_excel em massa 01/04/2025
# Re-running the CNPJ lookup appends the 3-company result set (IFOOD,
# GOOGLE, MERCADO PAGO) again to the bottom of the table, and the blank
# contact fields for GOOGLE on the existing row 6 collapse down to empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- company records: A..J = cnpj, nome, telefone, email, logradouro, bairro, municipio, uf, cep, atividade_principal
$ifood = @("33.157.312/0001-62","IFOOD BENEFICIOS E SERVICOS LTDA.","(11) 3634-3360","juridico@ifood.com.br","AV DOS AUTONOMISTAS 1496","VILA YARA","OSASCO","SP","06.020-902","Emissão de vales-alimentação, vales-transporte e similares")
$google = @("06.947.283/0001-60","GOOGLE INTERNATIONAL LLC","","","1600 AMPHITHEATER PARKWAY","","EXTERIOR","EX","","Holdings de instituições não-financeiras")
$mercadoPago = @("10.573.521/0001-91","MERCADO PAGO INSTITUICAO DE PAGAMENTO LTDA","(11) 2121-1212","naoresponder@mercadolivre.com","1 AV AVENIDA DAS NACOES UNIDAS, 3003","BONFIM","OSASCO","SP","06.233-903","Atividades de cobranças e informações cadastrais")

function Write-CnpjRow($rowIndex, $values) {
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $values[$col - 1]
    }
}

# Row 6 (GOOGLE) previously stored empty placeholders in C6/D6/F6/I6 - the
# new export drops those blank fields, so clear them.
$ws.Cells.Item(6, 3).Value = ""
$ws.Cells.Item(6, 4).Value = ""
$ws.Cells.Item(6, 6).Value = ""
$ws.Cells.Item(6, 9).Value = ""

# Append the same 3-company lookup result two more times (rows 8-15),
# matching the bulk re-query described in the commit message.
Write-CnpjRow 8 $ifood
Write-CnpjRow 9 $google
Write-CnpjRow 10 $ifood
Write-CnpjRow 11 $google
Write-CnpjRow 12 $mercadoPago
Write-CnpjRow 13 $ifood
Write-CnpjRow 14 $google
Write-CnpjRow 15 $mercadoPago
